# Api: add quiz and answer request api
#
# The "safety" sheet gains a 4th column ("audio_name") mirroring the one
# already present on the "departure" sheet: each category / sub-category
# row gets the filename of its narration clip.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("safety")

# Widen the new audio_name column (matches the other asset-name columns).
$ws.Columns.Item(4).ColumnWidth = 45

# sa_001 "Your rights and safety"
$ws.Range("D2").Value = "right_and_your_safety.mp3"
# sa_001_1 "Violence against women"
$ws.Range("D3").Value = "violence_on_women.mp3"
# sa_001_2 "Sexual harassment"
$ws.Range("D4").Value = "harassment.mp3"
# sa_001_3 "Your health"
$ws.Range("D5").Value = "your_health.mp3"
# sa_001_4 "Understand exploitation and human trafficking"
$ws.Range("D6").Value = "understand_exploitation_and_human_traficking.mp3"
# sa_002 "Safety planning"
$ws.Range("D8").Value = "my_body_my_choice_safety_planning_tips.mp3"

# Match the cell formatting used elsewhere in the sheet for these rows by
# copying each donor cell's format onto its new neighbour in column D.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("A5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("A8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$excel.CutCopyMode = $false
